# "resaltado verde y letra amarilla"
# Add a new paragraph "Main6" right after the "Main5" block (i.e. right
# after the third paragraph following "Main5": Main5 / video blurb /
# "Lorem ipsum ... elit."), styled like its sibling "Main*" paragraphs:
# green highlight on the run plus yellow (FFFF00) font color on both the
# run and the paragraph mark.

$d = $word.ActiveDocument

# Locate the unique "Main5" paragraph so the insertion point is found by
# content rather than a brittle hard-coded paragraph index.
$find = $d.Content
$find.Find.Execute("Main5", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$main5Index = $find.Paragraphs.Item(1).Index

# The anchor paragraph is two paragraphs after "Main5":
#   Main5 -> "El vídeo proporciona..." -> "Lorem ipsum ... elit."
$anchorIndex = $main5Index + 2
$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Insert a brand-new paragraph right after it.
$anchorPara.Range.InsertParagraphAfter()

# Grab the freshly created (empty) paragraph and fill it in.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newRange = $newPara.Range
$newRange.Text = "Main6"

# Highlight green (matches the rest of the "Main*" runs).
$newRange.HighlightColorIndex = 4

# Yellow font color (wdColorYellow / RGB(255,255,0) = 0x00FFFF00 -> 65535)
# applied to the whole paragraph range so it lands on both the run rPr
# and the paragraph-mark rPr, matching the target markup.
$newRange.Font.Color = 65535
